# Changes done by Khubim
# Add the two "Task Done" entries that document the data-cleaning work.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Load the dataset and checked for Null values"
$ws.Range("B3").Value = "Renamed the columns and  Cleaned city and country "

# Leave the active cell/selection on B3, matching where the author left off.
$ws.Range("B3").Select()
